# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" table (rows 16-20) is resequenced: the GRISEL ALCALA
# ARZUZA record (previously last, row 20) now leads the table, followed by
# the two workers' 1801 periods, then the same two workers' 1806 periods.
# Only the data cells (C:G) change; row/column styles stay as they were.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16: GRISEL ALCALA ARZUZA / periodo 1610
$ws.Range("C16").Value = "1047384097"
$ws.Range("D16").Value = "GRISEL ALCALA ARZUZA"
$ws.Range("E16").Value = "1610"
$ws.Range("F16").Value = 24640
$ws.Range("G16").Value = 616000

# Row 17: VALERIA MARGARITA HERRERA ORTEGON / periodo 1801
$ws.Range("C17").Value = "1007229713"
$ws.Range("D17").Value = "VALERIA MARGARITA HERRERA ORTEGON"
$ws.Range("E17").Value = "1801"
$ws.Range("F17").Value = 29509
$ws.Range("G17").Value = 781242

# Row 18: MARIA ALEJANDRA HERRERA TORRES / periodo 1801
$ws.Range("C18").Value = "1143366337"
$ws.Range("D18").Value = "MARIA ALEJANDRA HERRERA TORRES"
$ws.Range("E18").Value = "1801"
$ws.Range("F18").Value = 29509
$ws.Range("G18").Value = 781242

# Row 19: VALERIA MARGARITA HERRERA ORTEGON / periodo 1806
$ws.Range("C19").Value = "1007229713"
$ws.Range("D19").Value = "VALERIA MARGARITA HERRERA ORTEGON"
$ws.Range("E19").Value = "1806"
$ws.Range("F19").Value = 31249
$ws.Range("G19").Value = 781242

# Row 20: MARIA ALEJANDRA HERRERA TORRES / periodo 1806
$ws.Range("C20").Value = "1143366337"
$ws.Range("D20").Value = "MARIA ALEJANDRA HERRERA TORRES"
$ws.Range("E20").Value = "1806"
$ws.Range("F20").Value = 31249
$ws.Range("G20").Value = 781242
